# Apply edits: add 4 new slides (R tutorial content), reorder, and
# add text to the existing 4th slide's content placeholder.

$p = $ppt.ActivePresentation

# --- 1. Update existing slide 4 content placeholder -----------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "Open interest V/S volume"

# --- 2. Add four new slides (Title and Content layout) ---------------------
# Created in this order so slide IDs are assigned sequentially: 260,261,262,263

# Slide 5 (id 260): "Read data into R"
$s5 = $p.Slides.Add(5, 16)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Read data into R"
$body5 = $s5.Shapes.Item(2).TextFrame.TextRange
$body5.Text = "read_csv(“file_name.csv”)`rread_excel(here(“folder_name`",“data.xlsx`"),`r`t`t`t`tsheet=“sheet_name`")`rreadRDS(here(“folder_name`",“file_name.RDS`"))"
$body5.Paragraphs(3,1).ParagraphFormat.Bullet.Visible = 0

# Slide 6 (id 261): "Plot using GGPLOT library"
$s6 = $p.Slides.Add(6, 16)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Plot using GGPLOT library"

# Slide 7 (id 262): "Regression using dummies"
$s7 = $p.Slides.Add(7, 16)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Regression using dummies"

# Slide 8 (id 263): "Stargazer package"
$s8 = $p.Slides.Add(8, 16)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Stargazer package"

# --- 3. Reorder: move "Stargazer package" (id 263) right after "Read data into R" --
$p.Slides.Item(8).MoveTo(6)

Write-Host "Final slide count: $($p.Slides.Count)"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    Write-Host "$i : id=$($sl.SlideID) title=$($sl.Shapes.Item(1).TextFrame.TextRange.Text)"
}
